$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vals = @(
    "14+58=",
    "66-30=",
    "16+80=",
    "13+49=",
    "90-9=",
    "20+23=",
    "19-13=",
    "83+6=",
    "53+30=",
    "57-40=",
    "30+37=",
    "27-23=",
    "31+64=",
    "0+15=",
    "3+17=",
    "85-72=",
    "10+59=",
    "99-62=",
    "53+21=",
    "72+14=",
    "3+0=",
    "37-29=",
    "12+3=",
    "72+23=",
    "94-94=",
    "54-12=",
    "15+61=",
    "57+17=",
    "51-0=",
    "70+12=",
    "38-7=",
    "10+47=",
    "23+43=",
    "27+32=",
    "77-71=",
    "74+20=",
    "89-45=",
    "6+23=",
    "11+60=",
    "13+21=",
    "33+53=",
    "87-16=",
    "57+5=",
    "6+1=",
    "76-46=",
    "64-9=",
    "56+16=",
    "44-5=",
    "19+40=",
    "84+13=",
    "47-0=",
    "29+32=",
    "15-13=",
    "22+60=",
    "34-23=",
    "67+24=",
    "99-77=",
    "39-11=",
    "5+39=",
    "79-17=",
    "0+51=",
    "28+19=",
    "43+55=",
    "29-15=",
    "21-11=",
    "28+45=",
    "78-78=",
    "29-22=",
    "66+28=",
    "62-39=",
    "47+6=",
    "26+30=",
    "71-47=",
    "15+80=",
    "51-23=",
    "53-41=",
    "99-27=",
    "82-70=",
    "49+1=",
    "76-51=",
    "3+72=",
    "75+0=",
    "69-47=",
    "74-35=",
    "51-34=",
    "51+13=",
    "51+37=",
    "13+73=",
    "88-0=",
    "34+15=",
    "74-42=",
    "63-32=",
    "96-68=",
    "58+21=",
    "18-0=",
    "39-17=",
    "64-64=",
    "40+54=",
    "76-24=",
    "84-33="
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.MoveEnd(1, -1) | Out-Null
        $rng.Text = $vals[$idx]
        $idx++
    }
}
